$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "409.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.647"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.742"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.142"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.60"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000217"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.141"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "448.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "91.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "34.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.115"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.169"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0503"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.91"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.997"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.137"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.314"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "141.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "108.08"
$ws.Range("D50").Style = "Normal"

$ws.Range("D2").Value = "61.983.84"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "3.416.88"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("E6").Value = "  -3.78%  "
$ws.Range("E7").Value = "  +9.70%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +8.03%  "
$ws.Range("E10").Value = "  +16.85%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("B12").Value = "ShibaInu"
$ws.Range("C12").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("E12").Value = "  +65.74%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("E14").Value = "  +7.24%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.953.69"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("E16").Value = "  +4.61%  "
$ws.Range("D17").Value = "3.437.64"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("E18").Value = "  +8.88%  "
$ws.Range("E19").Value = "  +5.14%  "
$ws.Range("D20").Value = "61.931.52"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("E21").Value = "  +42.84%  "
$ws.Range("E22").Value = "  +8.32%  "
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("E25").Value = "  +2.94%  "
$ws.Range("E26").Value = "  +15.96%  "
$ws.Range("E27").Value = "  +7.47%  "
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("E30").Value = "  -3.05%  "
$ws.Range("E31").Value = "  +6.02%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E33").Value = "  -2.47%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("E36").Value = "  +4.15%  "
$ws.Range("E37").Value = "  +4.95%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  +9.28%  "
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("E41").Value = "  -1.10%  "
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("E43").Value = "  +2.29%  "
$ws.Range("E44").Value = "  -0.40%  "
$ws.Range("E45").Value = "  +2.53%  "
$ws.Range("E46").Value = "  +7.83%  "
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("E48").Value = "  +4.38%  "
$ws.Range("D49").Value = "3.760.03"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("E50").Value = "  +29.63%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.109.41"
$ws.Range("E51").Value = "  -0.48%  "
